$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 40
$ws.Range("A40").Value = 111880462
$ws.Range("B40").Value = 88966
$ws.Range("E40").Value = 5754
$ws.Range("F40").Value = "Gultoppig fingersvamp"
$ws.Range("G40").Value = "Ramaria testaceoflava"
$ws.Range("H40").Value = "(Bres.) Corner"
$ws.Range("I40").Value = "'1"
$ws.Range("Q40").Value = 509970.2466718731
$ws.Range("R40").Value = 6753250.046013334
$ws.Range("AO40").Value = "Pinus sylvestris # vid tallar"
$ws.Range("AL40").Value = "vid tallar"

# Row 42
$ws.Range("A42").Value = 111880475
$ws.Range("I42").Value = "'2"
$ws.Range("Q42").Value = 509957.7514087428
$ws.Range("R42").Value = 6753362.853637428
$ws.Range("AJ42").Value = "gran"
$ws.Range("AK42").Value = "Picea abies"
$ws.Range("AO42").Value = "Picea abies"
$ws.Range("AL42").ClearContents()

# Row 43
$ws.Range("A43").Value = 111880591
$ws.Range("B43").Value = 90658
$ws.Range("E43").Value = 4361
$ws.Range("F43").Value = "Orange taggsvamp"
$ws.Range("G43").Value = "Hydnellum aurantiacum"
$ws.Range("H43").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("I43").Value = "'8"
$ws.Range("Q43").Value = 509822.1902239832
$ws.Range("R43").Value = 6753234.069152902
$ws.Range("AJ43").Value = "tall"
$ws.Range("AK43").Value = "Pinus sylvestris"
$ws.Range("AO43").Value = "Pinus sylvestris"

# Row 44
$ws.Range("A44").Value = 111880500
$ws.Range("I44").Value = "'4"
$ws.Range("Q44").Value = 509899.1991435916
$ws.Range("R44").Value = 6753571.34232254

# Row 45
$ws.Range("A45").Value = 111880484
$ws.Range("B45").Value = 90658
$ws.Range("E45").Value = 4361
$ws.Range("F45").Value = "Orange taggsvamp"
$ws.Range("G45").Value = "Hydnellum aurantiacum"
$ws.Range("H45").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("I45").Value = "'11"
$ws.Range("Q45").Value = 509900.7891887496
$ws.Range("R45").Value = 6753525.142772059

# Row 46
$ws.Range("A46").Value = 111880601
$ws.Range("B46").Value = 88966
$ws.Range("E46").Value = 5754
$ws.Range("F46").Value = "Gultoppig fingersvamp"
$ws.Range("G46").Value = "Ramaria testaceoflava"
$ws.Range("H46").Value = "(Bres.) Corner"
$ws.Range("I46").Value = "'4"
$ws.Range("Q46").Value = 509941.5744066621
$ws.Range("R46").Value = 6753224.672924293

# Row 47
$ws.Range("A47").Value = 111880562
$ws.Range("I47").Value = "'3"
$ws.Range("Q47").Value = 509657.7198006394
$ws.Range("R47").Value = 6753521.069647122

# Row 48
$ws.Range("A48").Value = 111880574
$ws.Range("I48").Value = "'2"
$ws.Range("Q48").Value = 509595.7160662179
$ws.Range("R48").Value = 6753391.52735021
